# Issue 724 - Update DevMan diagrams to match changes to design
# Adjust the "Logic" component diagram: shrink/retime the boxes feeding the
# AccountsLogic/CoursesLogic/EvaluationsLogic stack to make room for a new
# "GateKeeper" box above AccountsLogic, and nudge a couple of unrelated
# boxes/connectors down to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TeamEvalResult box: move down slightly ------------------------------
$rectTeamEval = $s.Shapes.Item("Rectangle 85")
$rectTeamEval.Top = 282.0

# --- Emails box: move down slightly (sub-point precision to match EMU) --
$rectEmails = $s.Shapes.Item("Rectangle 115")
$rectEmails.Top = 318.11717224121094

# --- Connector feeding into the Emails box: follow it down --------------
$connToEmails = $s.Shapes.Item("Straight Arrow Connector 47")
$connToEmails.Top = 330.05857849121094

# --- AccountsLogic: move down and shrink to make room for GateKeeper ----
$rectAccounts = $s.Shapes.Item("Rectangle 39")
$rectAccounts.Top = 144.0
$rectAccounts.Height = 36.0

# --- CoursesLogic: move down and shrink ----------------------------------
$rectCourses = $s.Shapes.Item("Rectangle 55")
$rectCourses.Top = 186.0
$rectCourses.Height = 30.0

# --- EvaluationsLogic: move down and shrink ------------------------------
$rectEvaluations = $s.Shapes.Item("Rectangle 57")
$rectEvaluations.Top = 222.0
$rectEvaluations.Height = 30.0

# --- Dashed connectors on the right: follow their boxes down ------------
$connEvaluations = $s.Shapes.Item("Straight Arrow Connector 67")
$connEvaluations.Top = 198.0

$connCourses = $s.Shapes.Item("Straight Arrow Connector 79")
$connCourses.Top = 240.0

$connAccounts = $s.Shapes.Item("Straight Arrow Connector 81")
$connAccounts.Top = 162.0

# --- New "GateKeeper" box, styled like the other Logic boxes ------------
$gateKeeper = $rectAccounts.Duplicate()
$gateKeeper.Name = "Rectangle 54"
$gateKeeper.Left = 432.0
$gateKeeper.Top = 102.0
$gateKeeper.Width = 132.0
$gateKeeper.Height = 36.0
$gateKeeper.TextFrame.TextRange.Text = "GateKeeper"

# --- New dashed connector out of GateKeeper, matching its siblings ------
$gateKeeperConn = $connAccounts.Duplicate()
$gateKeeperConn.Name = "Straight Arrow Connector 58"
$gateKeeperConn.Left = 564.0
$gateKeeperConn.Top = 120.0
$gateKeeperConn.Width = 60.0
$gateKeeperConn.Height = 0.0
$gateKeeperConn.ConnectorFormat.BeginConnect($gateKeeper, 3)
